$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "last updated" date stamp in G1 (header row) ---
$ws.Range("G1").Value = 44268
$ws.Range("G1").NumberFormat = "d-mmm"

# --- Fill in "QuotesList" (File/Class Name) for rows that were missing it ---
$ws.Range("F5").Value = "QuotesList"
$ws.Range("F6").Value = "QuotesList"
$ws.Range("F8").Value = "QuotesList"
$ws.Range("F10").Value = "QuotesList"
$ws.Range("F26").Value = "QuotesList"
$ws.Range("F28").Value = "QuotesList"
$ws.Range("F29").Value = "QuotesList"

# --- Fill in missing Function Name values ---
$ws.Range("E9").Value = "def showRandomCryptoInHTML( )"
$ws.Range("E16").Value = "def getRandomQuote( )"
$ws.Range("E19").Value = "def deleteQuote(id)"
$ws.Range("E26").Value = "def addQuote(quote, author,submitted_by)"
$ws.Range("E28").Value = "def getStudentWithLeastQuotes()"
$ws.Range("E29").Value = "def updateQuote(id,quote, author,submitted_by)"

# --- Pull request status notes (new column G) ---
$ws.Range("G13").Value = "Waiting Pull Request"
$ws.Range("G25").Value = "Pull Request is in"

# --- New assignment row ---
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "Jasthi"
$ws.Range("C30").Value = "Siva"
$ws.Range("D30").Value = "Return the list of Quotes which contains some swear words"
$ws.Range("E30").Value = "def getQuotesContainingBadWords()"
$ws.Range("F30").Value = "QuotesList"

# --- Column G sizing for the new "Pull Request" notes column ---
$ws.Columns("G").ColumnWidth = 21

# --- Update the view's active selection to the last edited cell ---
$ws.Range("E29").Select()
